# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows (Especial, Extra (doble especial), Primera, Segunda)
# for Comercializadora del Agro de Limarí - Chirimoya, dated 2021-11-04 (44504),
# ahead of the existing historical rows (which shift down from 64-66 to 68-70).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four blank rows above the current row 64, shifting old rows 64-66 down to 68-70.
$ws.Rows.Item(64).Insert()
$ws.Rows.Item(64).Insert()
$ws.Rows.Item(64).Insert()
$ws.Rows.Item(64).Insert()

$newRows = @(
    @{ Row=64; L="Especial";                 M=400; N=1900;  O=2000;  P=1950;  S=1950; },
    @{ Row=65; L="Extra (doble especial)";    M=240; N=2100;  O=2200;  P=2150;  S=2150; },
    @{ Row=66; L="Primera";                   M=400; N=1600;  O=1700;  P=1650;  S=1650; },
    @{ Row=67; L="Segunda";                   M=300; N=1100;  O=1200;  P=1150;  S=1150; }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44504
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value = 100107002
    $ws.Cells.Item($row, 10).Value = "Chirimoya"
    $ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/kilo (en caja de 15 kilos)"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 1
}
